$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 4-10 (A:F) after the edit:
# Row 4  -> CustomNBC (values formerly held by the "NBC" row)
# Row 5  -> Kraken2_0.0   (unchanged values, shifted down from old row 4)
# Row 6  -> Kraken2_0.05  (unchanged values, shifted down from old row 5)
# Row 7  -> Kraken2_0.1   (unchanged values, shifted down from old row 6)
# Row 8  -> MMSeqs2_100   (unchanged values, shifted down from old row 7)
# Row 9  -> MMSeqs2_97    (unchanged values, shifted down from old row 8)
# Row 10 -> Mothur        (unchanged values, shifted down from old row 9)
# Rows 11 (Qiime2) and 12 (VSEARCH) remain untouched.

$data = @(
    @("CustomNBC",    0.49, 0.72, 0.65, 0.79, 0.5600000000000001),
    @("Kraken2_0.0",  0.58, 0.74, 0.72, 0.76, 0.7),
    @("Kraken2_0.05", 0.54, 0.75, 0.68, 0.8100000000000001, 0.59),
    @("Kraken2_0.1",  0.48, 0.75, 0.63, 0.87, 0.5),
    @("MMSeqs2_100",  0.54, 0.82, 0.6899999999999999, 0.9399999999999999, 0.55),
    @("MMSeqs2_97",   0.61, 0.8100000000000001, 0.74, 0.87, 0.65),
    @("Mothur",       0.45, 0.64, 0.6, 0.67, 0.54)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
